# NEVADA_2019.xlsx edit script
# 1. Rename header columns (row 1) to short snake_case codes.
# 2. Title-case the Spanish connector words (de, del, la, los, las, el, y)
#    inside every state/municipality name in columns A and B.
# 3. Fix two 1-ULP floating point values in D13 / D641.
# 4. Drop the trailing footer rows (1389-1394) that held sample-size /
#    source / author notes, shrinking the sheet to A1:D1388.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header row -----------------------------------------------------
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- 2. Title-case connector words in columns A and B -------------------
$lastRow = 1388
for ($r = 2; $r -le $lastRow; $r++) {
    for ($c = 1; $c -le 2; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $cell.Value()
        if ($v -ne $null -and $v -ne "") {
            $new = $v -replace '\bde\b', 'De'
            $new = $new -replace '\bdel\b', 'Del'
            $new = $new -replace '\bla\b', 'La'
            $new = $new -replace '\blos\b', 'Los'
            $new = $new -replace '\blas\b', 'Las'
            $new = $new -replace '\bel\b', 'El'
            $new = $new -replace '\by\b', 'Y'
            $cell.Value = $new
        }
    }
}

# --- 3. Floating point precision fixes ----------------------------------
$ws.Cells.Item(13, 4).Value = 0.009721048182586643
$ws.Cells.Item(641, 4).Value = 0.009056877188745319

# --- 4. Remove footer rows 1389-1394 ------------------------------------
$ws.Rows("1389:1394").Delete()
